$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="26.160.06"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Formula = '="1.663.31"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Formula = '="217.81"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("D6").Formula = '="0.5249"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Formula = '="1.003"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Formula = '="0.2642"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("D9").Formula = '="0.06281"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("D10").Formula = '="20.71"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -3.95%  '
$ws.Range("D11").Formula = '="0.07745"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Formula = '="4.463"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Formula = '="1.611.26"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  -3.08%  '
$ws.Range("D14").Formula = '="1.890.09"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Formula = '="0.5467"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Formula = '="0.0₅8120"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").Formula = '="64.87"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Formula = '="26.187.10"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").Formula = '="1.002"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Formula = '="4.588"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -2.82%  '
$ws.Range("D21").Formula = '="191.85"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").Formula = '="10.02"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -2.29%  '
$ws.Range("D23").Formula = '="6.003"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -4.03%  '
$ws.Range("D24").Formula = '="1.004"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("D25").Formula = '="137.65"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D26").Formula = '="0.1239"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -1.81%  '
$ws.Range("D27").Formula = '="7.255"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -1.63%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Formula = '="1.400"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("D30").Formula = '="0.05973"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -1.99%  '
$ws.Range("D31").Formula = '="1.280"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("D32").Formula = '="3.533"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -1.14%  '
$ws.Range("D33").Formula = '="3.265"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -3.59%  '
$ws.Range("D34").Formula = '="1.575"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -5.88%  '
$ws.Range("D35").Formula = '="0.9598"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  -3.73%  '
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").Formula = '="2.769"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="0.5668"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -5.93%  '
$ws.Range("D39").Formula = '="0.01598"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("D40").Formula = '="5.935"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("D41").Formula = '="0.8513"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("D43").Formula = '="101.09"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("D44").Formula = '="1.004.35"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  -7.22%  '
$ws.Range("D45").Formula = '="1.806.12"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("D46").Formula = '="56.73"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -1.76%  '
$ws.Range("E47").Value = '  -4.96%  '
$ws.Range("D48").Formula = '="0.9977"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("D49").Formula = '="7.999"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("D50").Formula = '="0.4315"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").Formula = '="0.05155"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -0.95%  '

$excel.CutCopyMode = $false
